$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "Neue Features" text in E26 (row 26) with the revised wording.
$ws.Range("E26").Value = "Wachen können im Jagdmodus nicht mehr zwischendrin in den Suchmodus wechseln; Wachen berechnen im Patrol- und Suchmodus nur noch 1x den Weg (im Jagdmodus pro Frame manuell neuen Path anfordern); Component für Pathfinding wird nicht komplett deaktiviert (sonst Error, wenn Seeker noch aktiv, stattdessen nur canMove auf false setzen); Erreichen des Ziels von Wachen (auch toter, die geschoben werden) führt nicht mehr zum Gewinnen; Tastaturbefehle werden jetzt immer erkannt, auch wenn Grafikeinstellungen auf `"Niedrigst`" sind (insb. auch Loslassen der Shift-Taste) (von FixedUpdate zu Update verschoben); Werfen des Steins ruft Wachen zur tatsächlich getroffenen Position (z.B. Wand, die getroffen wird), anstatt zur errechneten maximalen Poition gemäß Wurfreichweite; tote Wachen können nicht mehr alerted werden oder den Spieler jagen"

# Add the new row 27 with the new demo build entry.
$ws.Range("B27").Value = "DiscordiaAgency_Demo_2017_09_22-2.exe"
$ws.Range("C27").Value = "Entwicklung"
$ws.Range("D27").Value = "Anna Franziska"
$ws.Range("E27").Value = "Kugeln spawnen korrekt aus den Pistolen der Wachen anstatt aus Wachen-Mittelpunkt; Tutorial-Levels haben nur Features aktiviert, die auch benutzt werden sollen; SPACE führt jetzt immer weiter/zum nächsten Level, ENTER ist dazu da, das aktuelle Level (sofern gewonnen) zu wiederholen"

# Match styles of the row above (B26:E26) for the new row's cells.
$ws.Range("B26:D26").Copy()
$ws.Range("B27:D27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E26").Copy()
$ws.Range("E27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row heights.
$ws.Rows.Item(26).RowHeight = 255
$ws.Rows.Item(27).RowHeight = 90

# Update selection to mirror the authored state.
$ws.Range("E27").Select()
